$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Zapallo italiano" at
# "Vega Monumental Concepción". It belongs above the current row 178,
# so insert a fresh row there (this pushes the old rows 178-210 down to
# 179-211, growing the used range to A1:R211) and fill in its values.
$ws.Rows(178).Insert()

$ws.Range("A178").Value = 11
$ws.Range("B178").Value = "Vega Monumental Concepción"
$ws.Range("C178").Value = "Bíobío"
$ws.Range("D178").Value = 45015
$ws.Range("E178").Value = 8
$ws.Range("F178").Value = 100112032
$ws.Range("G178").Value = "Zapallo italiano"
$ws.Range("H178").Value = "Huracán"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 220
$ws.Range("K178").Value = 6000
$ws.Range("L178").Value = 6500
$ws.Range("M178").Value = 6273
$ws.Range("N178").Value = "$/caja 50 unidades"
$ws.Range("O178").Value = "Región de Arica y Parinacota"
$ws.Range("P178").Value = 125
$ws.Range("Q178").Value = 50
$ws.Range("R178").Value = "Hortaliza"
